$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "Encuesta_etapa_usuario"

# --- Column widths ---
$colWidths = @(8.71, 13.71, 11.71, 7.71, 13.71, 15.71, 17.71, 19.71, 9.71, 10.71, 12.71, 22.71, 8.71, 8.71, 8.71, 8.71, 8.71, 8.71)
for ($i = 0; $i -lt $colWidths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i]
}

# --- View: freeze header row, hide gridlines ---
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Header row style (row 1) ---
$header = $ws.Range("A1:R1")
$header.Font.Name = "Calibri"
$header.Font.Size = 11
$header.Font.Bold = $true
$header.Font.Color = 16777215
$header.Interior.Pattern = 1
$header.Interior.Color = 7949855
$header.Borders.Item(9).LineStyle = 1
$header.Borders.Item(9).Color = 0
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# --- Body styles (rows 2-23) ---
# NOTE: multi-area ("A,B,C") ranges only apply formatting to their first
# area in this runtime, so every contiguous block is handled individually.
$allBody = $ws.Range("A2:R23")
$allBody.Font.Name = "Calibri"
$allBody.Font.Size = 11
$allBody.Font.Bold = $false
$allBody.Font.Color = 0
$allBody.Borders.LineStyle = 1
$allBody.Borders.Color = 0
$allBody.VerticalAlignment = -4108

$leftAlignBlocks = @("A2:C23")
$numberBlocks = @("D2:E23", "G2:G23", "I2:J23", "L2:L23")
$percentBlocks = @("F2:F23", "H2:H23", "K2:K23")
$centerBlocks = @("M2:R23")

foreach ($addr in $leftAlignBlocks) {
    $ws.Range($addr).HorizontalAlignment = -4131
}
foreach ($addr in $numberBlocks) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4108
    $rng.NumberFormat = "#,##0"
}
foreach ($addr in $percentBlocks) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4108
    $rng.NumberFormat = "0.0""%"""
}
foreach ($addr in $centerBlocks) {
    $ws.Range($addr).HorizontalAlignment = -4108
}

# --- AutoFilter ---
$ws.Range("A1:R23").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Encuesta_etapa_usuario'!`$A`$1:`$R`$23")
$filterName.Visible = $false

Write-Host "Done"
